$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (rows 2-20): convert from text ("1") to numeric "disponible" counts,
# applying a numeric (thousand-separator) format, centered horizontally, top-aligned.
$fRange = $ws.Range("F2:F20")
$fRange.NumberFormat = "#,##0"
$fRange.HorizontalAlignment = -4108
$fRange.VerticalAlignment = -4160

$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 2.8
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 468
$ws.Range("F9").Value = 315
$ws.Range("F10").Value = 90
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 1.9999999999999876
$ws.Range("F16").Value = 3.0000000000000133
$ws.Range("F17").Value = -0.000000000000012434497875801753
$ws.Range("F18").Value = 5.0000000000000213
$ws.Range("F19").Value = 250
$ws.Range("F20").Value = 0

# E20: was a text reference to "43123020"; becomes a plain number.
$ws.Range("E20").Value = 43123020
